$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 141
$ws.Range("I2").Value = 325
$ws.Range("J2").Value = 1348
$ws.Range("K2").Value = 14
$ws.Range("L2").Value = 405
$ws.Range("M2").Value = 26
$ws.Range("N2").Value = 242
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 20
$ws.Range("S2").Value = 150
$ws.Range("T2").Value = 208
$ws.Range("U2").Value = 17
$ws.Range("V2").Value = 2087
$ws.Range("X2").Value = 2105
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 33
$ws.Range("AA2").Value = 18
